# feat: add 2022-Q1 data
#
# The workbook tracks quarterly fund-holding snapshots, one sheet per
# quarter, plus a "总计" (totals) summary sheet. This adds a new
# "2022-Q1" detail sheet (positioned right before "总计") and updates
# the "总计" summary sheet with a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet immediately before "总计".
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row (row 1), columns B..H.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1Header = $q1.Range("B1:H1")
$q1Header.Font.Bold = $true
$q1Header.HorizontalAlignment = -4108
$q1Header.VerticalAlignment = -4160
$q1Header.Borders.Item(7).LineStyle = 1
$q1Header.Borders.Item(8).LineStyle = 1
$q1Header.Borders.Item(9).LineStyle = 1
$q1Header.Borders.Item(10).LineStyle = 1

# Data rows 2..4 - fund holdings for 2022-Q1.
# Columns B,D,E,F,G hold numeric-looking text (fund codes / percentages)
# that must stay text (leading zeros, fixed decimals) - pre-format as
# Text so the values aren't coerced to numbers.
$q1.Range("B2:G4").NumberFormat = "@"

$q1Data = @(
    @("002692", "富国创新科技混合A",   "40.49", "92.07", "2.77", "1.1216", 10),
    @("011120", "富国创新科技混合C",   "0.90",  "92.07", "2.77", "0.0249", 10),
    @("519097", "新华中小市值优选混合", "0.75",  "62.70", "3.25", "0.0244", 7)
)

for ($i = 0; $i -lt $q1Data.Count; $i++) {
    $row = $i + 2
    $rec = $q1Data[$i]

    $aCell = $q1.Cells.Item($row, 1)
    $aCell.Value = $i
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.Item(7).LineStyle = 1
    $aCell.Borders.Item(8).LineStyle = 1
    $aCell.Borders.Item(9).LineStyle = 1
    $aCell.Borders.Item(10).LineStyle = 1

    $q1.Cells.Item($row, 2).Value = $rec[0]
    $q1.Cells.Item($row, 3).Value = $rec[1]
    $q1.Cells.Item($row, 4).Value = $rec[2]
    $q1.Cells.Item($row, 5).Value = $rec[3]
    $q1.Cells.Item($row, 6).Value = $rec[4]
    $q1.Cells.Item($row, 7).Value = $rec[5]
    $q1.Cells.Item($row, 8).Value = $rec[6]
}

# ------------------------------------------------------------------
# 2) Update "总计" with a new leading row for 2022-Q1, pushing the
#    existing rows down by one.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Read existing rows (2..5) before we overwrite them. Value2 (unlike
# Value) returns the real scalar through this COM bridge.
$existing = @()
for ($r = 2; $r -le 5; $r++) {
    $existing += , @(
        $total.Cells.Item($r, 2).Value2,
        $total.Cells.Item($r, 3).Value2,
        $total.Cells.Item($r, 4).Value2
    )
}

# Shift rows 2..5 down to 3..6.
$total.Range("B3:B6").NumberFormat = "@"
for ($r = 5; $r -ge 2; $r--) {
    $src = $existing[$r - 2]
    $dstRow = $r + 1

    $aCell = $total.Cells.Item($dstRow, 1)
    $aCell.Value = $dstRow - 2
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.Item(7).LineStyle = 1
    $aCell.Borders.Item(8).LineStyle = 1
    $aCell.Borders.Item(9).LineStyle = 1
    $aCell.Borders.Item(10).LineStyle = 1

    $total.Cells.Item($dstRow, 2).Value = $src[0]
    $total.Cells.Item($dstRow, 3).Value = $src[1]
    $total.Cells.Item($dstRow, 4).Value = $src[2]
}

# New row 2: 2022-Q1 totals.
$aCell2 = $total.Cells.Item(2, 1)
$aCell2.Value = 0
$aCell2.Font.Bold = $true
$aCell2.HorizontalAlignment = -4108
$aCell2.VerticalAlignment = -4160
$aCell2.Borders.Item(7).LineStyle = 1
$aCell2.Borders.Item(8).LineStyle = 1
$aCell2.Borders.Item(9).LineStyle = 1
$aCell2.Borders.Item(10).LineStyle = 1

$total.Range("B2").NumberFormat = "@"
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 1.17

Write-Host "2022-Q1 sheet added and summary updated"
